$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Test 14 -> Test 20 / User 14 -> User 20 / email -> htest20@...
$ws.Range("A2").Value = "Test 20"
$ws.Range("B2").Value = " User 20"
$ws.Range("C2").Value = "htest20@hspheretest.com"

# Row 3: Test 15 -> Test 21 / User 15 -> User 21 / email -> htest21@...
$ws.Range("A3").Value = "Test 21"
$ws.Range("B3").Value = "User 21"
$ws.Range("C3").Value = "htest21@hspheretest.com"

# Country column (N): IN -> US for both rows
$ws.Range("N2").Value = "US"
$ws.Range("N3").Value = "US"

# N2/N3 also pick up the style used elsewhere in the same column (e.g. O2/O3)
$ws.Range("O2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("O3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# T3: hsphere cluster id 1.0 -> 2.0
$ws.Range("T3").Value = 2.0
